# For every row in the "Recorded By" column (column G), reorder the
# comma-separated list of recorders so that an entry of "System"/"system"
# that currently appears first is moved to the end (i.e. the list is
# reversed). Rows that are a single value, or whose first entry is not
# "System"/"system", are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$col = 7  # column G = "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $text = $cell.Text

    if ($text -eq $null -or $text -eq "") {
        continue
    }

    $rawParts = $text -split ","
    if ($rawParts.Count -lt 2) {
        continue
    }

    $first = $rawParts[0].Trim()
    if ($first.ToLower() -ne "system") {
        continue
    }

    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    $count = $parts.Count
    $reversed = @()
    for ($i = $count - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i]
    }

    $newValue = $reversed -join ", "
    $cell.Value = $newValue
}
